$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.472.90"
$ws.Range("E2").Value = "  +2.06%  "
$ws.Range("D3").Value = "3.350.30"
$ws.Range("E3").Value = "  +3.14%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "192.81"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +4.48%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "594.78"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +2.08%  "
$ws.Range("E7").Value = "  +0.00%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.607"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.92%  "
$ws.Range("E9").Value = "  +2.75%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "6.73"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +1.68%  "
$ws.Range("E11").Value = "  +2.29%  "
$ws.Range("D12").Value = "3.932.39"
$ws.Range("E12").Value = "  +3.27%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.138"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +0.71%  "
$ws.Range("E14").Value = "  +1.74%  "
$ws.Range("D15").Value = "69.463.00"
$ws.Range("E15").Value = "  +2.09%  "
$ws.Range("E16").Value = "  +0.99%  "
$ws.Range("D17").Value = "3.345.97"
$ws.Range("E17").Value = "  +3.30%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "5.85"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.64%  "
$ws.Range("E19").Value = "  +2.08%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "428.48"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +8.04%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "7.74"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +1.84%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "73.25"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +2.70%  "
$ws.Range("E23").Value = "  +0.21%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "0.520"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.68%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.0000121"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +2.14%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.191"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +2.46%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "9.62"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("E28").Value = "  +0.93%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.01"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +2.11%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "5.67"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +1.06%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "23.11"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +1.14%  "
$ws.Range("E32").Value = "  +1.69%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "7.05"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.33%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.52"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +2.05%  "
$ws.Range("B36").Value = "Monero"
$ws.Range("C36").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "164.47"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +1.64%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.93"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +1.21%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "27.14"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +1.76%  "
$ws.Range("B39").Value = "Mantle"
$ws.Range("C39").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.811"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.38%  "
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "4.60"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.05%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "6.48"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.37%  "
$ws.Range("D42").Value = "2.746.74"
$ws.Range("E42").Value = "  +5.21%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "2.52"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +1.23%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "25.57"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +1.63%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0689"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.78%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "41.19"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.05%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "343.56"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +2.14%  "
$ws.Range("E48").Value = "  +1.32%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "32.54"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +4.55%  "
$ws.Range("E50").Value = "  +3.45%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "6.32"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.06%  "
